# This script re-applies a row-content rotation to worksheet "Artfynd".
# The underlying data rows (3-20, except 11 and 17 which are untouched)
# had their entire row contents rotated among each other, e.g. the values
# that used to live in row 4 now live in row 3, row 6's old values now
# live in row 4, etc. Column headers (row 1) and row 2 are untouched.
#
# The rotation is expressed as a set of cycles. For a cycle [r0, r1, r2, ..., rn]
# the NEW content of r0 is the OLD content of r1, the NEW content of r1 is the
# OLD content of r2, ..., and the NEW content of rn is the OLD content of r0.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastCol = "AY"
$cycles = @(
    @(3, 4, 6, 5),
    @(7, 14, 12, 9, 8),
    @(10, 13),
    @(15, 20, 16, 19, 18)
)

# A handful of columns hold date- or time-formatted text (e.g. "2023-09-24",
# "18:10") stored as plain strings (inlineStr) in the workbook. If we write
# such a string back through Value2 while the cell's number format is
# General, Excel will silently reinterpret it as a real date/time serial
# number, which would not match the source data (still plain text). To avoid
# that, mark those particular columns as Text before writing, then clear the
# formatting again afterwards so no stray style survives in the output.
$dateLikeCols = @("Y", "Z", "AA", "AB")

# Pre-compute the full set of rows participating in any cycle.
$allRows = @()
foreach ($cycle in $cycles) {
    foreach ($r in $cycle) {
        $allRows += $r
    }
}

# Snapshot the "before" values of every row involved, as well as mark the
# date-like columns as text so that the upcoming writes keep them as text.
$snapshot = @{}
foreach ($r in $allRows) {
    $rowRange = $ws.Range("A${r}:${lastCol}${r}")
    $snapshot[$r] = $rowRange.Value2
}

foreach ($r in $allRows) {
    foreach ($col in $dateLikeCols) {
        $ws.Range("${col}${r}").NumberFormat = "@"
    }
}

# Apply the rotation: new content of cycle[i] = old content of cycle[i+1]
# (wrapping around).
foreach ($cycle in $cycles) {
    $n = $cycle.Length
    for ($i = 0; $i -lt $n; $i++) {
        $target = $cycle[$i]
        $source = $cycle[($i + 1) % $n]
        $destRange = $ws.Range("A${target}:${lastCol}${target}")
        $destRange.Value2 = $snapshot[$source]
    }
}

# Remove the temporary text-number-format marking so the saved file does not
# carry extra style information that was not present before.
foreach ($r in $allRows) {
    foreach ($col in $dateLikeCols) {
        $ws.Range("${col}${r}").ClearFormats()
    }
}
